# Actualización automática 2025-07-04 17:15:07
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M18").Value = 1368.58
$ws1.Range("M22").Value = "4 de 20"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F18").Value = 1368.58
$ws2.Range("F22").Value = 14863.96

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D16").Value = 11909.82
$ws3.Range("E16").Value = 32356.42
$ws3.Range("F16").Value = 0.2690497318046439
$ws3.Range("D19").Value = 14863.96
$ws3.Range("E19").Value = 50514.03762291768
$ws3.Range("F19").Value = 0.2273541640986198
